$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.595.63"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "3.013.94"
$ws.Range("E3").Value = "  -0.56%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'586.36"
$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("D6").Value = "'147.74"
$ws.Range("E6").Value = "  -2.23%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "'0.527"
$ws.Range("E8").Value = "  -1.71%  "

$ws.Range("D9").Value = "3.013.48"
$ws.Range("E9").Value = "  -0.55%  "

$ws.Range("E10").Value = "  -1.84%  "

$ws.Range("D11").Value = "'5.82"
$ws.Range("E11").Value = "  +0.76%  "

$ws.Range("D12").Value = "'0.464"
$ws.Range("E12").Value = "  +3.66%  "

$ws.Range("D13").Value = "'0.0000231"
$ws.Range("E13").Value = "  -1.01%  "

$ws.Range("D14").Value = "'34.79"
$ws.Range("E14").Value = "  -3.78%  "

$ws.Range("E15").Value = "  +2.51%  "

$ws.Range("D16").Value = "3.517.07"
$ws.Range("E16").Value = "  -0.41%  "

$ws.Range("D17").Value = "'7.11"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").Value = "62.560.97"
$ws.Range("E18").Value = "  -0.62%  "

$ws.Range("D19").Value = "3.013.10"
$ws.Range("E19").Value = "  -0.83%  "

$ws.Range("D20").Value = "'460.37"
$ws.Range("E20").Value = "  -3.55%  "

$ws.Range("E21").Value = "  -0.96%  "

$ws.Range("D22").Value = "'0.692"
$ws.Range("E22").Value = "  -1.39%  "

$ws.Range("D23").Value = "'7.45"
$ws.Range("E23").Value = "  -0.39%  "

$ws.Range("D24").Value = "'81.79"
$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("E25").Value = "  -7.11%  "

$ws.Range("D26").Value = "'12.34"
$ws.Range("E26").Value = "  -2.62%  "

$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.20%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'10.01"
$ws.Range("E28").Value = "  -5.53%  "

$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.22%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'2.66"
$ws.Range("E30").Value = "  -0.31%  "

$ws.Range("D31").Value = "'7.03"
$ws.Range("E31").Value = "  -4.36%  "

$ws.Range("D32").Value = "'2.11"
$ws.Range("E32").Value = "  -3.65%  "

$ws.Range("D33").Value = "'28.08"
$ws.Range("E33").Value = "  +2.42%  "

$ws.Range("E34").Value = "  -0.63%  "

$ws.Range("D35").Value = "0.0₃0819"
$ws.Range("E35").Value = "  +0.95%  "

$ws.Range("E36").Value = "  -2.77%  "

$ws.Range("D37").Value = "'5.78"
$ws.Range("E37").Value = "  -1.97%  "

$ws.Range("E38").Value = "  -3.72%  "

$ws.Range("D39").Value = "'50.40"
$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("D40").Value = "'9.17"
$ws.Range("E40").Value = "  -0.56%  "

$ws.Range("D41").Value = "'2.92"
$ws.Range("E41").Value = "  -9.90%  "

$ws.Range("E42").Value = "  +9.12%  "

$ws.Range("D43").Value = "'392.86"
$ws.Range("E43").Value = "  -9.01%  "

$ws.Range("D44").Value = "'0.0360"
$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").Value = "'0.270"
$ws.Range("E45").Value = "  -5.52%  "

$ws.Range("D46").Value = "2.746.27"
$ws.Range("E46").Value = "  -2.63%  "

$ws.Range("D47").Value = "'37.54"
$ws.Range("E47").Value = "  -1.97%  "

$ws.Range("D48").Value = "'129.41"
$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.109"
$ws.Range("E50").Value = "  +0.27%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'2.20"
$ws.Range("E51").Value = "  -0.09%  "
